$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Cell="B2"; Value=0.2007042253521127},
    @{Cell="C2"; Value=0.5422535211267606},
    @{Cell="J2"; Value=0.007042253521126761},
    @{Cell="P2"; Value=0.1654929577464789},
    @{Cell="S2"; Value=0.08450704225352113},
    @{Cell="C3"; Value=0.01875},
    @{Cell="J3"; Value=0.025},
    @{Cell="P3"; Value=0.6875},
    @{Cell="S3"; Value=0.26875},
    @{Cell="J4"; Value=0.05714285714285714},
    @{Cell="O4"; Value=0.02857142857142857},
    @{Cell="P4"; Value=0.6857142857142857},
    @{Cell="S4"; Value=0.2285714285714286},
    @{Cell="B6"; Value=0.04522613065326633},
    @{Cell="D6"; Value=0.02010050251256281},
    @{Cell="F6"; Value=0.03517587939698492},
    @{Cell="J6"; Value=0.271356783919598},
    @{Cell="O6"; Value=0.01507537688442211},
    @{Cell="Q6"; Value=0.1306532663316583},
    @{Cell="R6"; Value=0.08040201005025126},
    @{Cell="S6"; Value=0.4020100502512563},
    @{Cell="B7"; Value=0.09523809523809523},
    @{Cell="D7"; Value=0.02164502164502164},
    @{Cell="F7"; Value=0.02164502164502164},
    @{Cell="J7"; Value=0.1818181818181818},
    @{Cell="O7"; Value=0.02164502164502164},
    @{Cell="Q7"; Value=0.1298701298701299},
    @{Cell="R7"; Value=0.08658008658008658},
    @{Cell="S7"; Value=0.4415584415584415},
    @{Cell="B8"; Value=0.08849557522123894},
    @{Cell="D8"; Value=0.01106194690265487},
    @{Cell="F8"; Value=0.04424778761061947},
    @{Cell="J8"; Value=0.1349557522123894},
    @{Cell="Q8"; Value=0.168141592920354},
    @{Cell="R8"; Value=0.07743362831858407},
    @{Cell="S8"; Value=0.4756637168141593},
    @{Cell="B9"; Value=0.07772020725388601},
    @{Cell="D9"; Value=0.0155440414507772},
    @{Cell="F9"; Value=0.08808290155440414},
    @{Cell="J9"; Value=0.1191709844559585},
    @{Cell="O9"; Value=0.01036269430051814},
    @{Cell="Q9"; Value=0.1865284974093264},
    @{Cell="R9"; Value=0.1036269430051813},
    @{Cell="S9"; Value=0.3989637305699482},
    @{Cell="B10"; Value=0.1111111111111111},
    @{Cell="D10"; Value=0.01610305958132045},
    @{Cell="E10"; Value=0.0008051529790660225},
    @{Cell="F10"; Value=0.07165861513687601},
    @{Cell="J10"; Value=0.1159420289855072},
    @{Cell="O10"; Value=0.01449275362318841},
    @{Cell="Q10"; Value=0.1819645732689211},
    @{Cell="R10"; Value=0.07971014492753623},
    @{Cell="S10"; Value=0.4082125603864734},
    @{Cell="G11"; Value=0.1450777202072539},
    @{Cell="J11"; Value=0.08549222797927461},
    @{Cell="K11"; Value=0.2098445595854922},
    @{Cell="L11"; Value=0.5440414507772021},
    @{Cell="S11"; Value=0.0155440414507772},
    @{Cell="G12"; Value=0.7},
    @{Cell="J12"; Value=0.2136363636363636},
    @{Cell="K12"; Value=0.00909090909090909},
    @{Cell="L12"; Value=0.03636363636363636},
    @{Cell="S12"; Value=0.04090909090909091},
    @{Cell="G13"; Value=0.675},
    @{Cell="J13"; Value=0.325},
    @{Cell="F15"; Value=0.01886792452830189},
    @{Cell="H15"; Value=0.1415094339622641},
    @{Cell="I15"; Value=0.09433962264150944},
    @{Cell="J15"; Value=0.3726415094339622},
    @{Cell="K15"; Value=0.07547169811320754},
    @{Cell="M15"; Value=0.01415094339622642},
    @{Cell="O15"; Value=0.04716981132075472},
    @{Cell="S15"; Value=0.2358490566037736},
    @{Cell="F16"; Value=0.005681818181818182},
    @{Cell="H16"; Value=0.1875},
    @{Cell="I16"; Value=0.06818181818181818},
    @{Cell="J16"; Value=0.3863636363636364},
    @{Cell="K16"; Value=0.1022727272727273},
    @{Cell="M16"; Value=0.005681818181818182},
    @{Cell="O16"; Value=0.03409090909090909},
    @{Cell="S16"; Value=0.2102272727272727},
    @{Cell="F17"; Value=0.01272264631043257},
    @{Cell="H17"; Value=0.1628498727735369},
    @{Cell="I17"; Value=0.08651399491094147},
    @{Cell="J17"; Value=0.3587786259541985},
    @{Cell="K17"; Value=0.1348600508905853},
    @{Cell="M17"; Value=0.02544529262086514},
    @{Cell="O17"; Value=0.09923664122137404},
    @{Cell="S17"; Value=0.1195928753180662},
    @{Cell="F18"; Value=0.02127659574468085},
    @{Cell="H18"; Value=0.2021276595744681},
    @{Cell="I18"; Value=0.07446808510638298},
    @{Cell="J18"; Value=0.4202127659574468},
    @{Cell="K18"; Value=0.06914893617021277},
    @{Cell="M18"; Value=0.005319148936170213},
    @{Cell="N18"; Value=0.005319148936170213},
    @{Cell="O18"; Value=0.05851063829787234},
    @{Cell="S18"; Value=0.1436170212765958},
    @{Cell="F19"; Value=0.0157819225251076},
    @{Cell="H19"; Value=0.2087517934002869},
    @{Cell="I19"; Value=0.08249641319942611},
    @{Cell="J19"; Value=0.3385939741750359},
    @{Cell="K19"; Value=0.1427546628407461},
    @{Cell="M19"; Value=0.01865136298421808},
    @{Cell="N19"; Value=0.002869440459110474},
    @{Cell="O19"; Value=0.06384505021520803},
    @{Cell="S19"; Value=0.1262553802008608}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$wb.Save()
